# The paragraph originally reads "...the conventional cycloconverter using
# uses 2 separate converters..." - a duplicated word ("using"/"uses"). The
# fix removes the stray "using " so the sentence reads "...cycloconverter
# uses 2 separate converters...". Word also relocates its automatic
# "_GoBack" (last-edit-location) bookmark to sit at the point of the edit,
# which splits the surrounding run there.

$d = $word.ActiveDocument

# --- 1. Delete the stray duplicate word "using " -------------------------
$find = $d.Content.Find
$find.Execute("using uses 2 separate converters") | Out-Null
$usingStart = $find.Parent.Start
$toDelete = $d.Range($usingStart, $usingStart + 6)
$toDelete.Text = ""

# --- 2. Re-establish the run boundaries that existed around this sentence -
# (the editing engine coalesces same-formatted runs it touches, so the
# boundaries that must stay distinct are restored with a zero-length
# bookmark add/delete, which forces a split without leaving markup behind)
$script:splitCounter = 0
function Split-RunAt($pos) {
    $script:splitCounter = $script:splitCounter + 1
    $tmpName = "TmpSplitMarker" + $script:splitCounter
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($tmpName, $r) | Out-Null
    $d.Bookmarks($tmpName).Delete()
}

$findFig = $d.Content.Find
$findFig.Execute("as represented in Fig. 1.") | Out-Null
Split-RunAt $findFig.Parent.End

$findAsCan = $d.Content.Find
$findAsCan.Execute("As can be observed") | Out-Null
Split-RunAt $findAsCan.Parent.Start

$findInverter = $d.Content.Find
$findInverter.Execute("inverter. The below model") | Out-Null
Split-RunAt ($findInverter.Parent.Start + ("inverter").Length)

# --- 3. Move "_GoBack" to sit right after "...conventional cycloconverter "
$findUses = $d.Content.Find
$findUses.Execute("uses 2 separate converters") | Out-Null
$goBackPos = $findUses.Parent.Start
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
